$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "practiceTracker"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "crossfitStats"

$ws2.Range("A1").Value = "Skill"
$ws2.Range("B1").Value = "Stat"
$ws2.Range("C1").Value = "Date"

$ws2.Range("A2").Value = "Deadlift"
$ws2.Range("B2").Value = "225lbs"
$ws2.Range("C2").Value = 44927
$ws2.Range("C2").NumberFormat = "mm-dd-yy"

$ws2.Range("A3").Value = "Bench 1x3"
$ws2.Range("B3").Value = "125lbs"
$ws2.Range("C3").Value = 45020

$ws2.Range("C2").Copy()
$ws2.Range("C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("D9").Select()
